$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Atas%Jenny%coreGivesNoEmail%1,        Bandy%Kenneth%coreGivesNoEmail%1,        Bradin%Stuart A.%coreGivesNoEmail%1,        Cadwallender%Bruce A.%coreGivesNoEmail%1,        Cinti%Sandro K.%coreGivesNoEmail%1,        Collins%Curtis D.%coreGivesNoEmail%1,        Goldberg%Janet%coreGivesNoEmail%1,        Holmes%Jennifer G.%coreGivesNoEmail%1,        Kim%Christopher%coreGivesNoEmail%1,        Krupansky%Frank%coreGivesNoEmail%1,        Lozon%Marie M.%coreGivesNoEmail%1,        Rodgers%Phillip E.%coreGivesNoEmail%1,        Shlafer%Jean%coreGivesNoEmail%1,        Wagner%Deborah%coreGivesNoEmail%1,        Wilkerson%William M.%coreGivesNoEmail%1,        Wright%Carrie M.%coreGivesNoEmail%1]'
$ws.Range("E3").Value = '[ Douglas M.%Fleming%null%1,       Douglas M.%Fleming%null%1]'
$ws.Range("C4").Value = '"Community Health Workers’ Palliative Care Learning Needs and Training: Results from a Partnership between a US University and a Rural Community Organization in Mpumalanga Province, South Africa"'
$ws.Range("E4").Value = '[Cathy%Campbell%xref no email%1,  Marianne%Baernholdt%xref no email%1]'
$ws.Range("F4").Value = '10.1353/hpu.2016.0078'
$ws.Range("G4").Value = 'CROSSREF'

# H4 target text looks like a date ("2023-05-24"); force literal text via
# a leading apostrophe so Excel does not coerce it into a date serial,
# then clear the resulting number-format override so the cell keeps the
# workbook default style (matches original unstyled H4 cell).
$ws.Range("H4").Value = "'2023-05-24"
$ws.Range("H4").ClearFormats()
